# Update unit price (N) and recompute total amount (O = M * N) for rows 2-4.
# O3:O4 are written as a single range-formula assignment so the engine stores
# them as a shared formula group (matches the authored workbook's XML shape).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 400
$ws.Range("O2").Formula = "=M2*N2"

$ws.Range("N3").Value = 300
$ws.Range("N4").Value = 200
$ws.Range("O3:O4").Formula = "=M3*N3"

# New number formats: Won-currency for unit price / total, plain thousands for quantity.
$ws.Range("N2:O4").NumberFormat = "_-[$₩-412]* #,##0_-;\-[$₩-412]* #,##0_-;_-[$₩-412]* ""-""??_-;_-@_-"
$ws.Range("M2:M4").NumberFormat = "#,##0_ "

# Widen the total-amount column so the currency format is fully visible.
$ws.Columns("O").ColumnWidth = 16.43

# Move the selection, as left by the author after editing.
$ws.Range("J16").Select() | Out-Null
